$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.026854260769989355
$ws.Range("B1").Value = 0.026854260751806719
$ws.Range("A2").Value = 0.011250905918788662
$ws.Range("B2").Value = -0.011250905968868187
$ws.Range("A3").Value = -0.014911512719387792
$ws.Range("B3").Value = 0.014911512690406623
$ws.Range("A4").Value = -0.056300293751726294
$ws.Range("B4").Value = 0.056300293718163565
